# Vendors p0.xlsx -- add "Items Not Found" sheet; an unrecognized shopping
# request ("Jane" / "sdfsdfsdf" / 20) gets appended to the Shopping List and
# mirrored onto the new sheet; Inventory stock counts are refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Inventory: refresh the Stock column (D) values.
# ---------------------------------------------------------------------
$inv = $wb.Worksheets.Item("Inventory")
$inv.Range("D2").Value = 999
$inv.Range("D3").Value = 977
$inv.Range("D4").Value = 1000
$inv.Range("D5").Value = 980
$inv.Range("D6").Value = 9998
$inv.Range("D7").Value = 9980
$inv.Range("D8").Value = 50
$inv.Range("D9").Value = 95
$inv.Range("D10").Value = 99

# ---------------------------------------------------------------------
# 2) Shopping List: a client ("Jane") requested an item not in the
#    inventory datatable ("sdfsdfsdf"), so a new row is recorded at the
#    top of the list (pushing the existing rows down).
# ---------------------------------------------------------------------
$sl = $wb.Worksheets.Item("Shopping List")
$sl.Rows.Item(2).Insert()
$sl.Range("A2").Value = "Jane"
$sl.Range("B2").Value = "sdfsdfsdf"
$sl.Range("C2").Value = 20

# ---------------------------------------------------------------------
# 3) Expenses: updated total for Webb.
# ---------------------------------------------------------------------
$exp = $wb.Worksheets.Item("Expenses")
$exp.Range("B3").Value = 36.6

# ---------------------------------------------------------------------
# 4) New "Items Not Found" sheet (appended after "Expenses") -- this is
#    where unmatched / out-of-stock shopping requests get written.
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$notFound = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$notFound.Name = "Items Not Found"

$notFound.Range("A1").Value = "Client"
$notFound.Range("B1").Value = "Item "
$notFound.Range("C1").Value = "Unaccounted Quantity "

$notFound.Range("A2").Value = "Jane"
$notFound.Range("B2").Value = "sdfsdfsdf"
$notFound.Range("C2").Value = 20

$notFound.Columns.Item(1).ColumnWidth = 14.45
$notFound.Columns.Item(2).ColumnWidth = 16.31
$notFound.Columns.Item(3).ColumnWidth = 22.74

# ---------------------------------------------------------------------
# 5) Restore selections / active sheet to match the saved workbook state.
# ---------------------------------------------------------------------
$sl.Range("G7").Select() | Out-Null
$notFound.Range("A2").Select() | Out-Null

$inv.Activate() | Out-Null
$inv.Range("C32").Select() | Out-Null
